# Task #51400: New sample file added with roles
#
# Adds five role columns (admin, reviewer, publisher, editor, content-creator)
# to the user-bulk-upload-sample sheet, tweaks a couple of sample names
# (apostrophes), fills in Y/N role flags for the three sample rows, and adds
# a highlighted blank row further down the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Headers: new role columns after the existing firstname/lastname/email ---
$ws.Range("D1").Value = "admin"
$ws.Range("E1").Value = "reviewer"
$ws.Range("F1").Value = "publisher"
$ws.Range("G1").Value = "editor"
$ws.Range("H1").Value = "content-creator"

# --- Row 2 (TestUserFirstName / F'des / testuser@test.com) : all roles Y ---
$ws.Range("B2").Value = "F'des"
$ws.Range("D2:H2").Value = "Y"

# --- Row 3 (Christy / Fernandes / test12email@test.com) : all roles N ---
$ws.Range("D3:H3").Value = "N"

# --- Row 4 (Test / test's / testemail123@test.com) : admin Y, rest N ---
$ws.Range("B4").Value = "test's"
$ws.Range("D4").Value = "Y"
$ws.Range("E4:H4").Value = "N"

# --- New highlighted blank row further down the sheet ---
$ws.Range("B6:C6").Interior.Color = 16777215
$ws.Range("B6:C6").HorizontalAlignment = -4131

# --- Styling pass -----------------------------------------------------
# Columns A-C (plus the plain header cells D1/G1/H1) share the sheet's
# normal "Arial / theme text colour" look.
$ws.Range("A1:D1").Font.ThemeColor = 1
$ws.Range("G1:H1").Font.ThemeColor = 1
$ws.Range("A2:C4").Font.ThemeColor = 1

# The role flag cells that are visually called out use a white highlight
# fill with the text left aligned.
$ws.Range("E1:F1").Interior.Color = 16777215
$ws.Range("E1:F1").HorizontalAlignment = -4131

$ws.Range("E3:H3").Interior.Color = 16777215
$ws.Range("E3:H3").HorizontalAlignment = -4131

$ws.Range("D4:H4").Interior.Color = 16777215
$ws.Range("D4:H4").HorizontalAlignment = -4131

Write-Host "Applied role-column sample update"
